$d = $word.ActiveDocument

# Locate the target paragraph precisely by its distinctive text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*emploi est parfait pour ceux qui aiment les ordinateurs*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Target paragraph not found"
}

$newParaXml = '<w:p w14:paraId="28B92795" w14:textId="77777777" w:rsidR="00AB4FE7" w:rsidRPr="005905F9" w:rsidRDefault="00AB4FE7" w:rsidP="00AB4FE7"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr></w:pPr><w:r w:rsidRPr="005905F9"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr><w:t>2. L' + [char]0x2019 + 'emploi est parfait pour ceux qui aiment les ordinateurs</w:t></w:r><w:r w:rsidR="004F126C"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr><w:t xml:space="preserve"> et les technologies qui changent le monde</w:t></w:r><w:r w:rsidRPr="005905F9"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + "`r`n" + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" mc:Ignorable="w14"><w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.Range.InsertXML($packageXml)
Write-Host "Done"
